$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (shared string) renames ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Updated GDP (column C) values and Colony (column AL) corrections ---
$ws.Range("C2").Value = 2771.04675450926
$ws.Range("C4").Value = 2934.187009790061
$ws.Range("C5").Value = 2870.311589353206
$ws.Range("C6").Value = 1873.394108966653
$ws.Range("C7").Value = 1460.056109840828
$ws.Range("AL8").Value = 1
$ws.Range("C8").Value = 5191.140356354663
$ws.Range("C9").Value = 1909.084588129339
$ws.Range("C10").Value = 6128.19547247793
$ws.Range("C11").Value = 4547.50930098406
$ws.Range("C12").Value = 4729.735976516416
$ws.Range("C13").Value = 1268.249210347625
$ws.Range("C14").Value = 3587.883798243964
$ws.Range("C15").Value = 471.181692645893
$ws.Range("C17").Value = 1286.515571617672
$ws.Range("C18").Value = 2839.92516805933
$ws.Range("C20").Value = 647.8358464534491
$ws.Range("C21").Value = 557.6321326261259
$ws.Range("C22").Value = 710.2742021758368
$ws.Range("C23").Value = 473.2998774917226
$ws.Range("C24").Value = 341.5541149051794
$ws.Range("C25").Value = 612.3436990512633
$ws.Range("C26").Value = 815.8736791314819
$ws.Range("C27").Value = 2983.242707849043
$ws.Range("C28").Value = 2898.942214704482
$ws.Range("C29").Value = 1904.346464968814
$ws.Range("C30").Value = 1503.870423231357
$ws.Range("AL31").Value = 1
$ws.Range("C31").Value = 5555.389721901988
$ws.Range("C32").Value = 1955.461557360978
$ws.Range("C33").Value = 6336.709213679884
$ws.Range("C34").Value = 4633.590358399045
$ws.Range("C35").Value = 5082.354756663512
$ws.Range("C36").Value = 1357.563719132622
$ws.Range("C37").Value = 3579.960081455846
$ws.Range("C38").Value = 492.3430015592067
$ws.Range("C40").Value = 1303.425880277445
$ws.Range("C41").Value = 2948.84548976845
$ws.Range("C43").Value = 670.2645481663891
$ws.Range("C44").Value = 587.1403830380606
$ws.Range("C45").Value = 711.1128122770988
$ws.Range("C46").Value = 1629.435089125503
$ws.Range("C47").Value = 369.2024078290272
$ws.Range("C48").Value = 644.763840173281
$ws.Range("C49").Value = 864.5379000312432
$ws.Range("C50").Value = 3083.80337578809
$ws.Range("C51").Value = 2965.153206179127
$ws.Range("C52").Value = 1939.33862702996
$ws.Range("C53").Value = 1577.487171555845
$ws.Range("C54").Value = 2024.117324382548
$ws.Range("C55").Value = 4921.848409120176
$ws.Range("C56").Value = 5360.226632400601
$ws.Range("C57").Value = 1410.426304742003
$ws.Range("C58").Value = 3156.723844635973
$ws.Range("C59").Value = 2999.422762626143
$ws.Range("C60").Value = 1982.009737844954
$ws.Range("C61").Value = 1657.651524528445
$ws.Range("C62").Value = 2094.024217383061
$ws.Range("C63").Value = 5122.180090208862
$ws.Range("C64").Value = 5642.578115155247
$ws.Range("C65").Value = 1469.177610078392
$ws.Range("C66").Value = 2379.668184479739
$ws.Range("C67").Value = 3748.449444923865
$ws.Range("C68").Value = 513.7390871590731
$ws.Range("C70").Value = 1325.930225429421
$ws.Range("C71").Value = 692.4450379203138
$ws.Range("C72").Value = 561.9866500362131
$ws.Range("C73").Value = 684.6474015015979
$ws.Range("C74").Value = 1652.714170143609
$ws.Range("C75").Value = 389.9389667216314
$ws.Range("C76").Value = 1000.829216794104
$ws.Range("C77").Value = 683.460336640684
$ws.Range("C78").Value = 869.6014949562591
$ws.Range("C80").Value = 2995.45235738661
$ws.Range("C81").Value = 3212.740625904757
$ws.Range("C82").Value = 3056.152683606517
$ws.Range("C83").Value = 2000.792448761861
$ws.Range("C84").Value = 1716.389195271215
$ws.Range("C85").Value = 2201.396847776877
$ws.Range("C86").Value = 5295.682695961288
$ws.Range("C87").Value = 5919.20956823756
$ws.Range("C88").Value = 1544.619247249133
$ws.Range("C89").Value = 2497.68592515536
$ws.Range("C90").Value = 3796.882621798447
$ws.Range("C91").Value = 534.5063430177229
$ws.Range("C93").Value = 1360.10887014004
$ws.Range("C94").Value = 711.0361291687414
$ws.Range("C95").Value = 565.0479699255185
$ws.Range("C96").Value = 680.3923729568069
$ws.Range("C97").Value = 1671.292192516047
$ws.Range("C98").Value = 419.1838602515346
$ws.Range("C99").Value = 1032.277326842402
$ws.Range("C100").Value = 698.3833464078615
$ws.Range("C101").Value = 872.1235974568563
$ws.Range("C103").Value = 3087.12349650562
$ws.Range("C104").Value = 3843.198240901342
$ws.Range("C105").Value = 2286.013198234259
$ws.Range("C106").Value = 720.1523351943922
$ws.Range("C107").Value = 1401.753174264641
$ws.Range("C108").Value = 3008.669179463094
$ws.Range("C109").Value = 5412.131646018807
$ws.Range("C110").Value = 3252.634165082374
$ws.Range("C111").Value = 449.4203771491282
$ws.Range("C112").Value = 2612.856880840196
$ws.Range("C113").Value = 3137.260298393558
$ws.Range("C114").Value = 2025.814194788851
$ws.Range("C115").Value = 1640.18070024053
$ws.Range("C116").Value = 1060.095015975378
$ws.Range("C117").Value = 707.8672001573369
$ws.Range("C118").Value = 3125.07948072635
$ws.Range("C119").Value = 1692.460946584157
$ws.Range("C120").Value = 558.2093442539386
$ws.Range("C121").Value = 1775.027517189621
$ws.Range("C122").Value = 5996.49696468919
$ws.Range("C124").Value = 723.2321880005983
$ws.Range("C125").Value = 1338.716747746975
$ws.Range("C126").Value = 886.4370030633224
$ws.Range("C127").Value = 3748.320622951519
$ws.Range("C128").Value = 2361.056581219794
$ws.Range("C129").Value = 726.6520119370772
$ws.Range("C130").Value = 1441.783971398429
$ws.Range("C131").Value = 3012.536723186288
$ws.Range("C132").Value = 5330.539154475424
$ws.Range("C133").Value = 3314.741082534716
$ws.Range("C134").Value = 482.6390663355013
$ws.Range("C135").Value = 2735.187532014817
$ws.Range("C136").Value = 3210.869677115934
$ws.Range("C137").Value = 2067.29003376698
$ws.Range("C138").Value = 1751.664428859304
$ws.Range("C139").Value = 1093.134170274031
$ws.Range("C140").Value = 729.7808175407341
$ws.Range("C141").Value = 3222.05417836739
$ws.Range("C142").Value = 1732.587316450496
$ws.Range("C143").Value = 579.0880693780265
$ws.Range("C144").Value = 1836.014008604312
$ws.Range("C145").Value = 6114.227214287786
$ws.Range("C147").Value = 767.6026452352251
$ws.Range("C148").Value = 1384.519227335143
$ws.Range("C149").Value = 900.3889853519216
$ws.Range("C150").Value = 2425.561644739583
$ws.Range("C151").Value = 747.8284752776283
$ws.Range("C152").Value = 1469.192636109792
$ws.Range("C153").Value = 2854.757682901436
$ws.Range("C154").Value = 5176.058803160127
$ws.Range("C155").Value = 3382.563653843273
$ws.Range("C156").Value = 514.0573067519859
$ws.Range("C157").Value = 3242.636921959078
$ws.Range("C158").Value = 2111.193164269742
$ws.Range("C159").Value = 1875.732161108182
$ws.Range("C160").Value = 1129.713195979213
$ws.Range("C161").Value = 749.2194349876407
$ws.Range("C162").Value = 3212.81539531051
$ws.Range("C163").Value = 1705.033923663474
$ws.Range("C164").Value = 584.2111078769213
$ws.Range("C165").Value = 1895.214690888655
$ws.Range("C166").Value = 6262.368904654469
$ws.Range("C168").Value = 792.6212731169028
$ws.Range("C169").Value = 1431.756130822538
$ws.Range("C170").Value = 909.5979669529498
